# daily auto push: 2025-10-12 09:24 UTC
# Appends the new day's data row (row 94) to the bottom of the sheet,
# matching the existing A:D column layout (date, weekday, hour, ranking).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 94

# Leading apostrophe forces the date-looking string to stay plain text
# (same representation as the existing date cells in column A) instead of
# Excel auto-converting "2025/10/12" into a date serial number.
$ws.Cells.Item($row, 1).Value = "'2025/10/12"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = "日"
$ws.Cells.Item($row, 3).Value = 18
$ws.Cells.Item($row, 4).Value = 35
